# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 3 (pushing existing rows 3-10 down
# to 4-11), keeping row 2 (the most recent prior entry) in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3; this shifts rows 3..10 down to 4..11
# and copies the row-above formatting (so D3 keeps the date style used by
# the rest of column D).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Cells.Item(3, 1).Value = 12
$ws.Cells.Item(3, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(3, 3).Value = "Metropolitana"
$ws.Cells.Item(3, 4).Value = 44453
$ws.Cells.Item(3, 5).Value = 13
$ws.Cells.Item(3, 6).Value = 100112013
$ws.Cells.Item(3, 7).Value = "Alcachofa"
$ws.Cells.Item(3, 8).Value = "Española"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 50
$ws.Cells.Item(3, 11).Value = 12000
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 13).Value = 12000
$ws.Cells.Item(3, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(3, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 16).Value = 400
$ws.Cells.Item(3, 17).Value = 30
$ws.Cells.Item(3, 18).Value = "Hortaliza"
